$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.991.08'

$ws.Range('D3').Value = '3.319.09'
$ws.Range('E3').Value = '  +6.30%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.73'
$ws.Range('E5').Value = '  +1.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.51'
$ws.Range('E6').Value = '  +5.15%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '3.319.13'
$ws.Range('E8').Value = '  +6.59%  '

$ws.Range('E10').Value = '  +3.35%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.54'
$ws.Range('E11').Value = '  +4.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.475'
$ws.Range('E12').Value = '  +4.40%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('E13').Value = '  +1.63%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.97'
$ws.Range('E14').Value = '  +2.70%  '

$ws.Range('D15').Value = '3.864.75'
$ws.Range('E15').Value = '  +6.32%  '

$ws.Range('E16').Value = '  +0.88%  '

$ws.Range('D17').Value = '3.315.54'
$ws.Range('E17').Value = '  +6.31%  '

$ws.Range('D18').Value = '64.059.21'
$ws.Range('E18').Value = '  +1.60%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.92'
$ws.Range('E19').Value = '  +4.06%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.42'
$ws.Range('E20').Value = '  +2.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.31'
$ws.Range('E21').Value = '  +1.31%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.738'
$ws.Range('E22').Value = '  +6.17%  '

$ws.Range('E23').Value = '  +3.70%  '

$ws.Range('E24').Value = '  +6.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.85'
$ws.Range('E25').Value = '  -1.57%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.78'
$ws.Range('E27').Value = '  +2.32%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.42'
$ws.Range('E28').Value = '  +6.53%  '

$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.22'
$ws.Range('E30').Value = '  +3.85%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.18'
$ws.Range('E31').Value = '  +5.75%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '29.63'
$ws.Range('E32').Value = '  +11.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  -1.48%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  +2.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.10'
$ws.Range('E35').Value = '  +2.40%  '

$ws.Range('E36').Value = '  +3.91%  '

$ws.Range('D37').Value = '0.0₃0760'
$ws.Range('E37').Value = '  +7.94%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.13'
$ws.Range('E38').Value = '  +1.97%  '

$ws.Range('E39').Value = '  +4.90%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '434.36'
$ws.Range('E40').Value = '  +2.93%  '

$ws.Range('D41').Value = '3.063.99'
$ws.Range('E41').Value = '  +5.60%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.46'
$ws.Range('E42').Value = '  +3.22%  '

$ws.Range('E43').Value = '  +3.28%  '

$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('E45').Value = '  +2.20%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.22'
$ws.Range('E46').Value = '  +5.08%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.59'
$ws.Range('E47').Value = '  +4.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '36.20'
$ws.Range('E48').Value = '  +15.46%  '

$ws.Range('E50').Value = '  +2.80%  '

$ws.Range('E51').Value = '  +2.17%  '
